# cryptos.xlsx -- price/volume refresh (GitHub Actions scheduled update)
# Commit: "Updated cryptos list on Sun Jun  4 15:41:25 UTC 2023 with GitHub Actions"
#
# The sheet stores Price (col D) and Volume(1h) (col E) as plain text so that
# values such as "27.276.10" (thousands-grouped, no real decimal meaning) or
# "10.10" (trailing zero significant) round-trip exactly. Column D also holds
# single-decimal numbers (e.g. "307.78") that Excel's smart-parser would silently
# coerce into real numbers (losing the trailing zero / exact text) if assigned
# directly via Range.Value on a General-formatted cell. For those we stage the
# text in a scratch cell formatted as Text ("@"), then Copy / PasteSpecial
# (xlPasteValues) it across -- paste-values carries the text payload without
# carrying the scratch cell's number format, so the destination cell keeps its
# original (General) style untouched, exactly like the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$scratch = $ws.Range("G1")
$scratch.NumberFormat = "@"
$needsScratch = $false

$ws.Range("D2").Value = "27.276.10"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.909.32"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.34%  "
$scratch.Value = "307.78"
$scratch.Copy()
$ws.Range("D5").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.27%  "
$scratch.Value = "0.5374"
$scratch.Copy()
$ws.Range("D7").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E7").Value = "  +3.34%  "
$scratch.Value = "0.3817"
$scratch.Copy()
$ws.Range("D8").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E8").Value = "  +1.21%  "
$scratch.Value = "0.07297"
$scratch.Copy()
$ws.Range("D9").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E9").Value = "  +0.31%  "
$scratch.Value = "22.05"
$scratch.Copy()
$ws.Range("D10").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E10").Value = "  +4.12%  "
$scratch.Value = "0.9025"
$scratch.Copy()
$ws.Range("D11").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E11").Value = "  -0.08%  "
$scratch.Value = "0.08201"
$scratch.Copy()
$ws.Range("D12").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E12").Value = "  -0.86%  "
$scratch.Value = "95.79"
$scratch.Copy()
$ws.Range("D13").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E13").Value = "  -0.76%  "
$scratch.Value = "5.346"
$scratch.Copy()
$ws.Range("D14").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E14").Value = "  +1.26%  "
$scratch.Value = "1.005"
$scratch.Copy()
$ws.Range("D15").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E15").Value = "  +0.49%  "
$scratch.Value = "0.000008657"
$scratch.Copy()
$ws.Range("D16").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "27.310.87"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "1.144.44"
$ws.Range("E20").Value = "  -39.88%  "
$scratch.Value = "5.048"
$scratch.Copy()
$ws.Range("D21").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E21").Value = "  -0.80%  "
$scratch.Value = "10.82"
$scratch.Copy()
$ws.Range("D22").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E22").Value = "  +1.49%  "
$scratch.Value = "6.519"
$scratch.Copy()
$ws.Range("D23").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E23").Value = "  +1.42%  "
$scratch.Value = "149.96"
$scratch.Copy()
$ws.Range("D24").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E24").Value = "  +1.97%  "
$scratch.Value = "2.291"
$scratch.Copy()
$ws.Range("D25").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("E26").Value = "  +0.37%  "
$scratch.Value = "1.747"
$scratch.Copy()
$ws.Range("D27").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E27").Value = "  +0.03%  "
$scratch.Value = "117.02"
$scratch.Copy()
$ws.Range("D28").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("E29").Value = "  -0.11%  "
$scratch.Value = "4.808"
$scratch.Copy()
$ws.Range("D30").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E30").Value = "  -1.87%  "
$scratch.Value = "0.09301"
$scratch.Copy()
$ws.Range("D31").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E31").Value = "  +0.56%  "
$scratch.Value = "0.8389"
$scratch.Copy()
$ws.Range("D32").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E32").Value = "  +5.24%  "
$scratch.Value = "0.05070"
$scratch.Copy()
$ws.Range("D33").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  -0.86%  "
$scratch.Value = "3.007"
$scratch.Copy()
$ws.Range("D35").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E35").Value = "  +1.71%  "
$scratch.Value = "3.358"
$scratch.Copy()
$ws.Range("D36").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E36").Value = "  -1.92%  "
$scratch.Value = "2.705"
$scratch.Copy()
$ws.Range("D37").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E37").Value = "  +3.94%  "
$scratch.Value = "0.5761"
$scratch.Copy()
$ws.Range("D38").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E38").Value = "  +0.97%  "
$scratch.Value = "0.02009"
$scratch.Copy()
$ws.Range("D39").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("E40").Value = "  -0.08%  "
$scratch.Value = "9.314"
$scratch.Copy()
$ws.Range("D41").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E41").Value = "  +3.14%  "
$scratch.Value = "6.566"
$scratch.Copy()
$ws.Range("D42").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E42").Value = "  -0.36%  "
$scratch.Value = "117.56"
$scratch.Copy()
$ws.Range("D43").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("E44").Value = "  +0.50%  "
$scratch.Value = "0.4934"
$scratch.Copy()
$ws.Range("D45").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("E46").Value = "  +0.23%  "
$scratch.Value = "10.10"
$scratch.Copy()
$ws.Range("D47").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("E48").Value = "  +0.66%  "
$scratch.Value = "38.54"
$scratch.Copy()
$ws.Range("D49").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E49").Value = "  +2.41%  "
$scratch.Value = "0.06129"
$scratch.Copy()
$ws.Range("D50").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E50").Value = "  +2.97%  "
$scratch.Value = "63.45"
$scratch.Copy()
$ws.Range("D51").PasteSpecial($xlPasteValues)
$needsScratch = $true
$ws.Range("E51").Value = "  -0.79%  "

if ($needsScratch) {
    $scratch.Clear()
}
$ws.Range("A1").Select()
